$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.333402
$ws.Range("H2").Value = 73.00020600000001
$ws.Range("I2").Value = 0.9697175080062574
$ws.Range("J2").Value = 0.9697175080062576
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 33.380049
$ws.Range("N2").Value = 100.140147
$ws.Range("O2").Value = 0.3891462059670435
$ws.Range("P2").Value = 0.3891462059670435
$ws.Range("Q2").Value = 812.2501510966981
$ws.Range("R2").Value = 7310.251359870283
$ws.Range("S2").Value = 0.3773618891004513
$ws.Range("T2").Value = 0.3773618891004513

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.333402
$ws.Range("H3").Value = 73.00020600000001
$ws.Range("I3").Value = 0.9697175080062574
$ws.Range("J3").Value = 0.9697175080062576
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.85786133333334
$ws.Range("N3").Value = 101.573584
$ws.Range("O3").Value = 0.3947165649764305
$ws.Range("P3").Value = 0.3947165649764305
$ws.Range("Q3").Value = 823.8769506842563
$ws.Range("R3").Value = 7414.892556158306
$ws.Range("S3").Value = 0.3827635637577342
$ws.Range("T3").Value = 0.3827635637577342

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.333402
$ws.Range("H4").Value = 73.00020600000001
$ws.Range("I4").Value = 0.9697175080062574
$ws.Range("J4").Value = 0.9697175080062576
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.53974466666667
$ws.Range("N4").Value = 55.61923400000001
$ws.Range("O4").Value = 0.2161372290565261
$ws.Range("P4").Value = 0.2161372290565261
$ws.Range("Q4").Value = 451.1350599513561
$ws.Range("R4").Value = 4060.215539562205
$ws.Range("S4").Value = 0.2095920551480721
$ws.Range("T4").Value = 0.2095920551480721

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.288887
$ws.Range("H5").Value = 0.866661
$ws.Range("I5").Value = 0.01151252018667195
$ws.Range("J5").Value = 0.01151252018667195
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.380049
$ws.Range("N5").Value = 100.140147
$ws.Range("O5").Value = 0.3891462059670435
$ws.Range("P5").Value = 0.3891462059670435
$ws.Range("Q5").Value = 9.643062215463001
$ws.Range("R5").Value = 86.787559939167
$ws.Range("S5").Value = 0.00448005355176239
$ws.Range("T5").Value = 0.00448005355176239

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.288887
$ws.Range("H6").Value = 0.866661
$ws.Range("I6").Value = 0.01151252018667195
$ws.Range("J6").Value = 0.01151252018667195
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 33.85786133333334
$ws.Range("N6").Value = 101.573584
$ws.Range("O6").Value = 0.3947165649764305
$ws.Range("P6").Value = 0.3947165649764305
$ws.Range("Q6").Value = 9.781095987002669
$ws.Range("R6").Value = 88.02986388302401
$ws.Range("S6").Value = 0.004544182422304968
$ws.Range("T6").Value = 0.004544182422304968

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.288887
$ws.Range("H7").Value = 0.866661
$ws.Range("I7").Value = 0.01151252018667195
$ws.Range("J7").Value = 0.01151252018667195
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.53974466666667
$ws.Range("N7").Value = 55.61923400000001
$ws.Range("O7").Value = 0.2161372290565261
$ws.Range("P7").Value = 0.2161372290565261
$ws.Range("Q7").Value = 5.355891217519334
$ws.Range("R7").Value = 48.20302095767401
$ws.Range("S7").Value = 0.002488284212604596
$ws.Range("T7").Value = 0.002488284212604597

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4710003333333333
$ws.Range("H8").Value = 1.413001
$ws.Range("I8").Value = 0.01876997180707065
$ws.Range("J8").Value = 0.01876997180707065
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.380049
$ws.Range("N8").Value = 100.140147
$ws.Range("O8").Value = 0.3891462059670435
$ws.Range("P8").Value = 0.3891462059670435
$ws.Range("Q8").Value = 15.722014205683
$ws.Range("R8").Value = 141.498127851147
$ws.Range("S8").Value = 0.007304263314829914
$ws.Range("T8").Value = 0.007304263314829915

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4710003333333333
$ws.Range("H9").Value = 1.413001
$ws.Range("I9").Value = 0.01876997180707065
$ws.Range("J9").Value = 0.01876997180707065
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.85786133333334
$ws.Range("N9").Value = 101.573584
$ws.Range("O9").Value = 0.3947165649764305
$ws.Range("P9").Value = 0.3947165649764305
$ws.Range("Q9").Value = 15.94706397395378
$ws.Range("R9").Value = 143.523575765584
$ws.Range("S9").Value = 0.00740881879639137
$ws.Range("T9").Value = 0.007408818796391371

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4710003333333333
$ws.Range("H10").Value = 1.413001
$ws.Range("I10").Value = 0.01876997180707065
$ws.Range("J10").Value = 0.01876997180707065
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.53974466666667
$ws.Range("N10").Value = 55.61923400000001
$ws.Range("O10").Value = 0.2161372290565261
$ws.Range("P10").Value = 0.2161372290565261
$ws.Range("Q10").Value = 8.732225917914889
$ws.Range("R10").Value = 78.59003326123401
$ws.Range("S10").Value = 0.004056889695849365
$ws.Range("T10").Value = 0.004056889695849365
